$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D98").Value = 44581
$ws.Range("J98").Value = 230
$ws.Range("M98").Value = 9239
$ws.Range("D99").Value = 44449
$ws.Range("J99").Value = 250
$ws.Range("K99").Value = 9000
$ws.Range("L99").Value = 9500
$ws.Range("M99").Value = 9240
$ws.Range("P99").Value = 1540
$ws.Range("D100").Value = 44223
$ws.Range("J100").Value = 125
$ws.Range("L100").Value = 12000
$ws.Range("M100").Value = 11480
$ws.Range("P100").Value = 1913
$ws.Range("D101").Value = 44257
$ws.Range("J101").Value = 120
$ws.Range("K101").Value = 11000
$ws.Range("L101").Value = 11000
$ws.Range("M101").Value = 11000
$ws.Range("P101").Value = 1833
$ws.Range("I102").Value = 'Primera'
$ws.Range("J102").Value = 160
$ws.Range("K102").Value = 10000
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = 10000
$ws.Range("P102").Value = 1667
$ws.Range("D103").Value = 44188
$ws.Range("I103").Value = 'Segunda'
$ws.Range("J103").Value = 50
$ws.Range("K103").Value = 7000
$ws.Range("L103").Value = 7000
$ws.Range("M103").Value = 7000
$ws.Range("P103").Value = 1167
$ws.Range("I104").Value = 'Primera'
$ws.Range("J104").Value = 150
$ws.Range("K104").Value = 9000
$ws.Range("L104").Value = 9000
$ws.Range("M104").Value = 9000
$ws.Range("P104").Value = 1500
$ws.Range("D105").Value = 44189
$ws.Range("I105").Value = 'Segunda'
$ws.Range("J105").Value = 80
$ws.Range("K105").Value = 7000
$ws.Range("L105").Value = 7000
$ws.Range("M105").Value = 7000
$ws.Range("P105").Value = 1167
$ws.Range("D106").Value = 44321
$ws.Range("J106").Value = 160
$ws.Range("K106").Value = 9000
$ws.Range("L106").Value = 9000
$ws.Range("M106").Value = 9000
$ws.Range("P106").Value = 1500
$ws.Range("D107").Value = 44291
$ws.Range("J107").Value = 255
$ws.Range("K107").Value = 10000
$ws.Range("M107").Value = 10490
$ws.Range("P107").Value = 1748
$ws.Range("I108").Value = 'Primera'
$ws.Range("J108").Value = 110
$ws.Range("K108").Value = 11000
$ws.Range("L108").Value = 11000
$ws.Range("M108").Value = 11000
$ws.Range("P108").Value = 1833
$ws.Range("D109").Value = 44285
$ws.Range("I109").Value = 'Segunda'
$ws.Range("J109").Value = 60
$ws.Range("D110").Value = 44371
$ws.Range("J110").Value = 180
$ws.Range("D111").Value = 44397
$ws.Range("J111").Value = 160
$ws.Range("L111").Value = 9000
$ws.Range("M111").Value = 9000
$ws.Range("P111").Value = 1500
$ws.Range("D112").Value = 44336
$ws.Range("J112").Value = 110
$ws.Range("L112").Value = 9500
$ws.Range("M112").Value = 9227
$ws.Range("P112").Value = 1538
$ws.Range("D113").Value = 44342
$ws.Range("J113").Value = 160
$ws.Range("K113").Value = 9000
$ws.Range("L113").Value = 9000
$ws.Range("M113").Value = 9000
$ws.Range("P113").Value = 1500
$ws.Range("D114").Value = 44263
$ws.Range("J114").Value = 230
$ws.Range("L114").Value = 12000
$ws.Range("M114").Value = 11478
$ws.Range("P114").Value = 1913
$ws.Range("D115").Value = 44307
$ws.Range("J115").Value = 130
$ws.Range("K115").Value = 11000
$ws.Range("L115").Value = 11000
$ws.Range("M115").Value = 11000
$ws.Range("P115").Value = 1833
$ws.Range("I116").Value = 'Primera'
$ws.Range("J116").Value = 78
$ws.Range("K116").Value = 12000
$ws.Range("L116").Value = 12000
$ws.Range("M116").Value = 12000
$ws.Range("P116").Value = 2000
$ws.Range("D117").Value = 44166
$ws.Range("I117").Value = 'Segunda'
$ws.Range("J117").Value = 50
$ws.Range("K117").Value = 8500
$ws.Range("L117").Value = 8500
$ws.Range("M117").Value = 8500
$ws.Range("O117").Value = 'Pan de Azúcar'
$ws.Range("P117").Value = 1417
$ws.Range("D118").Value = 44442
$ws.Range("K118").Value = 9000
$ws.Range("L118").Value = 9000
$ws.Range("M118").Value = 9000
$ws.Range("O118").Value = 'Provincia del Elquí'
$ws.Range("P118").Value = 1500
$ws.Range("I119").Value = 'Primera'
$ws.Range("J119").Value = 120
$ws.Range("K119").Value = 12000
$ws.Range("L119").Value = 12000
$ws.Range("M119").Value = 12000
$ws.Range("P119").Value = 2000
$ws.Range("D120").Value = 44200
$ws.Range("I120").Value = 'Segunda'
$ws.Range("J120").Value = 90
$ws.Range("K120").Value = 8000
$ws.Range("L120").Value = 8000
$ws.Range("M120").Value = 8000
$ws.Range("O120").Value = 'Provincia de Santiago'
$ws.Range("P120").Value = 1333
$ws.Range("J121").Value = 730
$ws.Range("K121").Value = 8500
$ws.Range("L121").Value = 10000
$ws.Range("M121").Value = 9164
$ws.Range("O121").Value = 'Pan de Azúcar'
$ws.Range("P121").Value = 1527
$ws.Range("D122").Value = 44435
$ws.Range("J122").Value = 240
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 9500
$ws.Range("M122").Value = 9229
$ws.Range("O122").Value = 'Provincia de Limarí'
$ws.Range("P122").Value = 1538
$ws.Range("D123").Value = 44302
$ws.Range("J123").Value = 160
$ws.Range("K123").Value = 10000
$ws.Range("L123").Value = 10000
$ws.Range("M123").Value = 10000
$ws.Range("O123").Value = 'Provincia de Santiago'
$ws.Range("P123").Value = 1667
